# Timesheet update: add a new entry for row 8 (date + hours worked),
# which shifts the SUM/Total formulas in row 5 and moves the active
# cell selection down to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entry: date 10/25/2012 (serial 41207), 3 hours worked.
$ws.Range("A8").Value = 41207
$ws.Range("B8").Value = 3

# Recalculate dependent formulas (Total Hours / Invoice in row 5).
$wb.Application.Calculate()

# Move the current selection to the newly filled cell.
$ws.Range("B8").Select()
